$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.798.44'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.088.20'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = '  +1.23%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.72'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  -0.57%  '
$ws.Range("E6").Value = '  -0.18%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '57.35'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = '  +1.37%  '
$ws.Range("E9").Value = '  +1.56%  '
$ws.Range("E10").Value = '  +2.30%  '
$ws.Range("E11").Value = '  +3.03%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.383.85'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = '  +0.71%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.41'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = '  -1.15%  '
$ws.Range("E14").Value = '  +2.37%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.764'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = '  -1.32%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.23'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  +2.21%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.092.33'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = '  +1.40%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.668.89'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.13'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  -2.85%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.80'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  +1.95%  '
$ws.Range("E21").Value = '  +1.53%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '227.92'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  +0.80%  '
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("E24").Value = '  -1.55%  '
$ws.Range("E25").Value = '  -0.12%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '168.14'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  +1.17%  '
$ws.Range("E27").Value = '  +12.08%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.92'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = '  +1.95%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.44'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  -0.39%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.46'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  +2.26%  '
$ws.Range("E31").Value = '  +0.62%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.61'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  +3.73%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0624'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  +1.53%  '
$ws.Range("E34").Value = '  +0.28%  '
$ws.Range("E35").Value = '  -0.19%  '
$ws.Range("E36").Value = '  +3.44%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.999'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  -0.16%  '
$ws.Range("E39").Value = '  -4.70%  '
$ws.Range("E40").Value = '  +6.37%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.94'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = '  -0.44%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '96.96'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  +0.93%  '
$ws.Range("E43").Value = '  +0.76%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.452.86'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  -0.74%  '
$ws.Range("E45").Value = '  -0.26%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.06'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  +3.91%  '
$ws.Range("E47").Value = '  -3.58%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '15.61'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  +4.26%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.33'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = '  +2.70%  '
$ws.Range("E50").Value = '  +2.13%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.279.24'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  +1.15%  '
